# "Improved performance of exponential functions."
#
# Renames the "ops" summary-row label to "operators", refreshes the
# exp / exp2 / expm1 timing data (rows 8-10) and a handful of outlier
# cells in rows 22 ("asin"), 28 ("tanh") and 31 ("atanh") that had been
# manually flagged with the red "Bad" cell style - that highlight no
# longer applies, so it is removed along with the now-unused "Bad"
# cell style. All dependent ratio columns (P/Q/R, Z/AA/AB) and the
# "functions"/"all" summary rows (37/38) are plain formulas, so they
# recompute automatically once the raw timings below are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row label rename: "ops" -> "operators" ---
$ws.Range("A36").Value = "operators"

# --- exp (row 8) ---
$ws.Range("F8").Value  = 4689
$ws.Range("G8").Value  = 4338
$ws.Range("H8").Value  = 3788
$ws.Range("M8").Value  = 2105
$ws.Range("N8").Value  = 2178
$ws.Range("O8").Value  = 1192
$ws.Range("W8").Value  = 2733
$ws.Range("X8").Value  = 2222
$ws.Range("Y8").Value  = 1483

# --- exp2 (row 9) ---
$ws.Range("F9").Value  = 4481
$ws.Range("G9").Value  = 4688
$ws.Range("H9").Value  = 4074
$ws.Range("M9").Value  = 8026
$ws.Range("N9").Value  = 4416
$ws.Range("O9").Value  = 3160
$ws.Range("W9").Value  = 7239
$ws.Range("X9").Value  = 4336
$ws.Range("Y9").Value  = 3398

# --- expm1 (row 10) ---
$ws.Range("F10").Value = 5501
$ws.Range("G10").Value = 5254
$ws.Range("H10").Value = 4708
$ws.Range("M10").Value = 4766
$ws.Range("N10").Value = 2528
$ws.Range("O10").Value = 1467
$ws.Range("W10").Value = 4404
$ws.Range("X10").Value = 2627
$ws.Range("Y10").Value = 1680

# --- row 22 ("asin"): refresh values and drop the red "Bad" highlight ---
$ws.Range("H22").Value = 1172
$ws.Range("O22").Value = 301
$ws.Range("Y22").Value = 536
$ws.Range("H22").Style = "Normal"

# --- row 28 ("tanh"): refresh values and drop the red "Bad" highlight ---
$ws.Range("H28").Value = 2246
$ws.Range("O28").Value = 796
$ws.Range("Y28").Value = 1158
$ws.Range("H28").Style = "Normal"

# --- row 31 ("atanh"): refresh values and drop the red "Bad" highlight ---
$ws.Range("H31").Value = 1135
$ws.Range("Y31").Value = 781
$ws.Range("H31").Style = "Normal"

# The "Bad" cell style is no longer referenced anywhere - remove it.
$wb.Styles("Bad").Delete()

# Make sure every dependent formula (P/Q/R, Z/AA/AB ratios and the
# row 37/38 summaries) is refreshed against the new raw timings.
$excel.Calculate()
